# Auto-generated data-driven update of cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '57.412.57'; E = $null }
    @{ Row = 3; D = '3.011.56'; E = $null }
    @{ Row = 4; D = $null; E = '  -0.08%  ' }
    @{ Row = 5; D = '507.71'; E = '  -0.63%  ' }
    @{ Row = 6; D = '139.37'; E = '  +0.29%  ' }
    @{ Row = 7; D = $null; E = '  +0.05%  ' }
    @{ Row = 8; D = $null; E = '  -0.02%  ' }
    @{ Row = 9; D = '7.59'; E = '  +0.52%  ' }
    @{ Row = 10; D = $null; E = '  +1.17%  ' }
    @{ Row = 11; D = '0.366'; E = '  +2.32%  ' }
    @{ Row = 12; D = '3.529.08'; E = '  +0.13%  ' }
    @{ Row = 13; D = $null; E = '  +0.32%  ' }
    @{ Row = 14; D = '26.36'; E = '  +2.16%  ' }
    @{ Row = 15; D = $null; E = '  +3.64%  ' }
    @{ Row = 16; D = '57.408.11'; E = '  +1.10%  ' }
    @{ Row = 17; D = '6.19'; E = '  +3.82%  ' }
    @{ Row = 18; D = '3.014.69'; E = '  +0.20%  ' }
    @{ Row = 19; D = '12.82'; E = '  +2.30%  ' }
    @{ Row = 20; D = '7.94'; E = '  +0.98%  ' }
    @{ Row = 21; D = '327.29'; E = '  -1.15%  ' }
    @{ Row = 22; D = $null; E = '  -0.23%  ' }
    @{ Row = 23; D = $null; E = '  -1.88%  ' }
    @{ Row = 24; D = $null; E = '  +3.31%  ' }
    @{ Row = 25; D = '64.41'; E = '  +2.12%  ' }
    @{ Row = 26; D = $null; E = '  -3.58%  ' }
    @{ Row = 27; D = '0.999'; E = '  -0.05%  ' }
    @{ Row = 28; D = '0.0₃0919'; E = '  +1.28%  ' }
    @{ Row = 29; D = '6.78'; E = '  +1.00%  ' }
    @{ Row = 30; D = $null; E = '  +3.36%  ' }
    @{ Row = 31; D = $null; E = '  +0.43%  ' }
    @{ Row = 32; D = $null; E = '  -5.74%  ' }
    @{ Row = 33; D = '20.57'; E = '  -0.65%  ' }
    @{ Row = 35; D = '153.66'; E = '  -0.50%  ' }
    @{ Row = 36; D = $null; E = '  +3.60%  ' }
    @{ Row = 37; D = $null; E = '  -0.05%  ' }
    @{ Row = 38; D = '24.54'; E = '  +2.24%  ' }
    @{ Row = 39; D = $null; E = '  -0.42%  ' }
    @{ Row = 40; D = '3.046.60'; E = '  +0.18%  ' }
    @{ Row = 41; D = '37.81'; E = '  +2.27%  ' }
    @{ Row = 42; D = '3.86'; E = '  +4.78%  ' }
    @{ Row = 43; D = $null; E = '  -0.04%  ' }
    @{ Row = 44; D = '0.649'; E = '  -0.13%  ' }
    @{ Row = 45; D = $null; E = '  -0.62%  ' }
    @{ Row = 46; D = '2.221.21'; E = '  -2.30%  ' }
    @{ Row = 47; D = '0.981'; E = '  -1.97%  ' }
    @{ Row = 48; D = $null; E = '  +3.71%  ' }
    @{ Row = 49; D = '0.0239'; E = '  -0.35%  ' }
    @{ Row = 50; D = '19.49'; E = '  -0.79%  ' }
    @{ Row = 51; D = $null; E = '  -6.04%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Range("D" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.ClearFormats()
    }
    if ($null -ne $u.E) {
        $cell = $ws.Range("E" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
        $cell.ClearFormats()
    }
}

Write-Output "Updated $($updates.Count) rows"
